$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the tiny floating point precision of the existing row 16 timestamp
$ws.Range("A16").Value = 45874.54183811343

# Append the new sensor reading as row 17
$ws.Range("A17").Value = 45874.5836159344
$ws.Range("A17").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B17").Value = 2025
$ws.Range("C17").Value = 19
$ws.Range("D17").Value = 20.56
$ws.Range("E17").Value = 74.25
$ws.Range("F17").Value = 82.62
$ws.Range("G17").Value = 13.87
$ws.Range("H17").Value = "ESE"
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = "14:00:24"
